# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker table (rows 16-21, columns B:G) is regrouped: instead of being
# grouped by period (1804 block then 1805 block), it is now grouped by
# worker, each worker showing period 1805 then 1804. The "Salario Basico"
# (column G) values for LINO RICARDO LEON BOLIVAR and CARLOS ANDRES LEON
# FRANCO are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: LINO RICARDO LEON BOLIVAR - periodo 1805
$ws.Range("C16").Value = "73145513"
$ws.Range("D16").Value = "LINO RICARDO LEON BOLIVAR"
$ws.Range("E16").Value = "1805"
$ws.Range("F16").Value = 29269
$ws.Range("G16").Value = 731717

# Row 17: LINO RICARDO LEON BOLIVAR - periodo 1804
$ws.Range("C17").Value = "73145513"
$ws.Range("D17").Value = "LINO RICARDO LEON BOLIVAR"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 29269
$ws.Range("G17").Value = 731717

# Row 18: CARLOS ANDRES LEON FRANCO - periodo 1805
$ws.Range("C18").Value = "1047457155"
$ws.Range("D18").Value = "CARLOS ANDRES LEON FRANCO"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 781242

# Row 19: CARLOS ANDRES LEON FRANCO - periodo 1804
$ws.Range("C19").Value = "1047457155"
$ws.Range("D19").Value = "CARLOS ANDRES LEON FRANCO"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

# Row 20: JORGE RAFAEL LEON FRANCO - periodo 1805
$ws.Range("C20").Value = "1047471570"
$ws.Range("D20").Value = "JORGE RAFAEL LEON FRANCO"
$ws.Range("E20").Value = "1805"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242

# Row 21: JORGE RAFAEL LEON FRANCO - periodo 1804
$ws.Range("C21").Value = "1047471570"
$ws.Range("D21").Value = "JORGE RAFAEL LEON FRANCO"
$ws.Range("E21").Value = "1804"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242
